$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.336.92"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.542.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +9.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +9.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "552.57"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.535.93"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +9.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.634"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +15.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.88"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.102.48"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.539.49"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +9.92%  "
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.424.10"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +7.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.31%  "
$ws.Range("E20").Value = "  +7.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +17.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.17%  "
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.13"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.12"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  +9.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.36"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "644.74"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.70"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.74"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("E34").Value = "  +4.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0831"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +15.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.68"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.69%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.147"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +18.94%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.33"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +13.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.041.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("E45").Value = "  +11.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.19%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.98%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0418"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.81%  "
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +11.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.05%  "
